$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-06 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.07521283979117906
$ws.Range("E2").Value = 0.002462271644162017
$ws.Range("D3").Value = 0.0462636746381683
$ws.Range("E3").Value = -0.0009018418026918162
$ws.Range("D4").Value = 0.03690110387357954
$ws.Range("E4").Value = -0.004858072027943861
$ws.Range("D5").Value = 0.03349883105267772
$ws.Range("E5").Value = -0.001780496712929125
$ws.Range("D6").Value = 0.03083566327855625
$ws.Range("E6").Value = -0.007030334591850096
$ws.Range("D7").Value = 0.03181463694672623
$ws.Range("E7").Value = -0.004371417240508935
$ws.Range("D8").Value = 0.02968052787728362
$ws.Range("E8").Value = -0.0002447531053051399
$ws.Range("D9").Value = 0.02887207647273551
$ws.Range("E9").Value = 0.001522881291403211
$ws.Range("D10").Value = 0.0263726444820133
$ws.Range("E10").Value = -0.0003696857670981712
$ws.Range("D11").Value = 0.0273891827519789
$ws.Range("E11").Value = 0.0003664345914256018
$ws.Range("D12").Value = 0.02342291519088601
$ws.Range("E12").Value = 0.006525198938992061
$ws.Range("D13").Value = 0.02396678945098044
$ws.Range("E13").Value = -0.002763819095477293
$ws.Range("D14").Value = 0.02027398836759056
$ws.Range("E14").Value = 0.005657548584198535
$ws.Range("D15").Value = 0.01966358748446876
$ws.Range("E15").Value = -0.01368850865253735
$ws.Range("D16").Value = 0.0207897130723954
$ws.Range("E16").Value = -0.007345926349933252
$ws.Range("D17").Value = 0.0186582325692362
$ws.Range("E17").Value = 0.004805278634440135
$ws.Range("D18").Value = 0.01778898959536998
$ws.Range("E18").Value = 0.003562447611064501
$ws.Range("D19").Value = 0.01485139950126061
$ws.Range("E19").Value = -0.002831858407079557
$ws.Range("D20").Value = 0.0136961685157595
$ws.Range("E20").Value = -0.007160354249105105
$ws.Range("D21").Value = 0.01623980821215723
$ws.Range("E21").Value = -0.008578550386844186
$ws.Range("D22").Value = 0.01408223321462091
$ws.Range("E22").Value = 0.00339378801042578
$ws.Range("D23").Value = 0.01309054684406208
$ws.Range("E23").Value = -0.003942958533219465
$ws.Range("D24").Value = 0.01497049744995616
$ws.Range("E24").Value = 0.008791924454575106
$ws.Range("D25").Value = 0.01415334788061568
$ws.Range("E25").Value = -0.01569506726457404
$ws.Range("D26").Value = 0.0125454299889973
$ws.Range("E26").Value = -0.002773333333333405
$ws.Range("D27").Value = 0.01223478124641261
$ws.Range("E27").Value = -0.001171875000000044
$ws.Range("D28").Value = 0.01240234039895137
$ws.Range("E28").Value = -0.01472798316801938
$ws.Range("D29").Value = 0.01152363458646487
$ws.Range("E29").Value = 0.0007299270072993469
$ws.Range("D30").Value = 0.01263241163692101
$ws.Range("E30").Value = -0.01259079903147686
$ws.Range("D31").Value = 0.01304151213484793
$ws.Range("E31").Value = -0.01759014951627069
$ws.Range("D32").Value = 0.01356626954924485
$ws.Range("E32").Value = -0.005136334812936072
$ws.Range("D33").Value = 0.01112098602396086
$ws.Range("E33").Value = 0.007198263821740936
$ws.Range("D34").Value = 0.01153543923868309
$ws.Range("E34").Value = 0.004076779344317938
$ws.Range("D35").Value = 0.009908022164453238
$ws.Range("E35").Value = 0.0008248317777297398
$ws.Range("D36").Value = 0.01085277667882466
$ws.Range("E36").Value = 0.007195606892633988
$ws.Range("D37").Value = 0.01062418699639657
$ws.Range("E37").Value = 0.003958614484930401
$ws.Range("D38").Value = 0.01020552807562089
$ws.Range("E38").Value = -0.008710311885360977
$ws.Range("D39").Value = 0.009338196786324422
$ws.Range("E39").Value = -0.009437438584998348
$ws.Range("D40").Value = 0.009165380501218843
$ws.Range("E40").Value = 0.0009177373602535788
$ws.Range("D41").Value = 0.009389334348565287
$ws.Range("E41").Value = -0.0161761562032352
$ws.Range("D42").Value = 0.009040738667271016
$ws.Range("E42").Value = 0.01509769094138536
$ws.Range("D43").Value = 0.009567981271608609
$ws.Range("E43").Value = -0.01571428571428557
$ws.Range("D44").Value = 0.009923936938496424
$ws.Range("E44").Value = 0.001618122977346426
$ws.Range("D45").Value = 0.008990556947315025
$ws.Range("E45").Value = -0.008898669983733565
$ws.Range("D46").Value = 0.009363048685731198
$ws.Range("E46").Value = -0.006339581036383879
$ws.Range("D47").Value = 0.008770044132194135
$ws.Range("E47").Value = 0.002615746795710194
$ws.Range("D48").Value = 0.007069027202028833
$ws.Range("E48").Value = -0.01130401860565733
$ws.Range("D49").Value = 0.008191377212930213
$ws.Range("E49").Value = 0.004562533548040815
$ws.Range("D50").Value = 0.007958820785019888
$ws.Range("E50").Value = 0.01918573230048648
$ws.Range("D51").Value = 0.007818646513942824
$ws.Range("E51").Value = -0.01313593770056909
$ws.Range("D52").Value = 0.00753127253099486
$ws.Range("E52").Value = -0.005597014925373234
$ws.Range("D53").Value = 0.007204661222252736
$ws.Range("E53").Value = 0.02786069651741285
$ws.Range("D54").Value = 0.007493994681884696
$ws.Range("E54").Value = -0.002946353409351787
$ws.Range("D55").Value = 0.006657680474732454
$ws.Range("E55").Value = 0.0008646495100677054
$ws.Range("D56").Value = 0.006574283235377023
$ws.Range("E56").Value = 0.005931956964233764
$ws.Range("D57").Value = 0.006609075894546509
$ws.Range("E57").Value = -0.003471017007983357
$ws.Range("D58").Value = 0.006258377360225511
$ws.Range("E58").Value = -0.002061855670103085
$ws.Range("D59").Value = 0.00547845784788095
$ws.Range("E59").Value = 0.005565684675175264
$ws.Range("D60").Value = 0.006641192195318342
$ws.Range("E60").Value = 0.007987910189982683
$ws.Range("D61").Value = 0.005312571419340723
$ws.Range("E61").Value = 0.003688377114069663
$ws.Range("D62").Value = 0.005756464576437131
$ws.Range("E62").Value = 0.005446333687566307
$ws.Range("D63").Value = 0.005292689899815304
$ws.Range("E63").Value = -0.002745069710322845
$ws.Range("D64").Value = 0.004829679897021375
$ws.Range("E64").Value = -0.0186827105763141
$ws.Range("D65").Value = 0.004685347712005103
$ws.Range("E65").Value = -0.001917662899343142
$ws.Range("D66").Value = 0.004371544689880316
$ws.Range("E66").Value = -0.01134798294522799
$ws.Range("D67").Value = 0.00427949707784676
$ws.Range("E67").Value = 0.01701956580005359
$ws.Range("D68").Value = 0.003599032955243546
$ws.Range("E68").Value = 0.008644729503624937
$ws.Range("D69").Value = 0.004037047682287966
$ws.Range("E69").Value = 0.002130908832617129
$ws.Range("D70").Value = 0.003637935736237999
$ws.Range("E70").Value = 0.00144508670520227
$ws.Range("D71").Value = 0.00315843745402996
$ws.Range("E71").Value = 0.003631576558173366
$ws.Range("D72").Value = 0.002673968791940567
$ws.Range("E72").Value = -0.009008042895442325
$ws.Range("D73").Value = 0.002583976240819493
$ws.Range("E73").Value = 0.007139290140011578
$ws.Range("D74").Value = 0.002349555920453658
$ws.Range("E74").Value = -0.0005695455839876962
$ws.Range("D75").Value = 0.001832636412792726
$ws.Range("E75").Value = -0.009075262087310132
$ws.Range("D76").Value = 0.001880810863950475
$ws.Range("E76").Value = 0.01743151903237306
$ws.Range("E77").Value = -0.001454029673358948

$ws.Protect()
